$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-07-06 Sunday" "2025-07-07 Monday"

Replace-Text "544×6=3264" "295×5=1475"
Replace-Text "505×6=3030" "440×8=3520"
Replace-Text "320×3=960" "819×6=4914"
Replace-Text "486×4=1944" "300×6=1800"
Replace-Text "774×4=3096" "998×3=2994"

Replace-Text "310×6=1860" "941×4=3764"
Replace-Text "291×5=1455" "673×2=1346"
Replace-Text "707×6=4242" "105×3=315"
Replace-Text "617×3=1851" "875×4=3500"
Replace-Text "942×2=1884" "777×6=4662"

Replace-Text "316×6=1896" "582×2=1164"
Replace-Text "107×5=535" "213×2=426"
Replace-Text "829×8=6632" "966×9=8694"
Replace-Text "254×6=1524" "378×9=3402"
Replace-Text "639×5=3195" "756×3=2268"

Replace-Text "176×2=352" "380×9=3420"
Replace-Text "710×2=1420" "925×4=3700"
Replace-Text "375×5=1875" "468×2=936"
Replace-Text "615×3=1845" "474×8=3792"
Replace-Text "861×9=7749" "703×7=4921"

Replace-Text "850×7=5950" "449×9=4041"
Replace-Text "488×7=3416" "354×2=708"
Replace-Text "350×4=1400" "821×3=2463"
Replace-Text "513×9=4617" "804×8=6432"
Replace-Text "494×2=988" "208×5=1040"
